$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1:3").Delete()
$ws.Rows(13).RowHeight = 36
$ws.Rows(14).RowHeight = 23.1
$ws.Rows(15).RowHeight = 21
Write-Output ("UsedRange rows: " + $ws.UsedRange.Rows.Count)
Write-Output ("Dimension: " + $ws.UsedRange.Address())
